$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data rows to append (April 2025 update without ECB)
$newData = @(
    @{ Row = 562; Date = 45536; Mean = 3.6300718468268096; Median = 3 },
    @{ Row = 563; Date = 45566; Mean = 3.4590725565951042; Median = 3 },
    @{ Row = 564; Date = 45597; Mean = 3.3850959326505774; Median = 3 },
    @{ Row = 565; Date = 45627; Mean = 2.9430802164335881; Median = 3 },
    @{ Row = 566; Date = 45658; Mean = 3.451973780215102;  Median = 3 },
    @{ Row = 567; Date = 45689; Mean = 3.7916655071943244; Median = 3 }
)

# Copy formatting from the last existing row (561) down to the new rows
$ws.Range("A561:C561").Copy()
$ws.Range("A562:C567").PasteSpecial(-4122)

foreach ($item in $newData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.Mean
    $ws.Cells.Item($r, 3).Value = $item.Median
}
